$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" footer timestamp text (A1)
$ws.Range("A1").Value = "Datos actualizados a 6 de Agosto de 2020 a las 23:52"

# Swap the Santa Lucia / Timor Oriental rows (A202 <-> A203 text contents)
$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("A203").Value = "Timor Oriental"

# Row 4 (row id 8) - Estados Unidos-ish stats update
$ws.Range("B4").Value = 5022591
$ws.Range("C4").Value = 49023
$ws.Range("D4").Value = 2568622
$ws.Range("E4").Value = 2291360
$ws.Range("G4").Value = 1008
$ws.Range("H4").Value = 162609

# Row 5
$ws.Range("B5").Value = 2912212
$ws.Range("C5").Value = 49451
$ws.Range("E5").Value = 793082
$ws.Range("G5").Value = 1075
$ws.Range("H5").Value = 98493

# Row 8
$ws.Range("B8").Value = 538184
$ws.Range("C8").Value = 8307
$ws.Range("D8").Value = 387316
$ws.Range("E8").Value = 141264
$ws.Range("G8").Value = 306
$ws.Range("H8").Value = 9604

# Row 23
$ws.Range("D23").Value = 82460
$ws.Range("E23").Value = 82861

# Row 30
$ws.Range("B30").Value = 95006
$ws.Range("C30").Value = 131
$ws.Range("D30").Value = 48898
$ws.Range("E30").Value = 41157
$ws.Range("G30").Value = 21
$ws.Range("H30").Value = 4951

# Row 79
$ws.Range("E79").Value = 6399
$ws.Range("G79").Value = 3
$ws.Range("H79").Value = 92

# Row 121
$ws.Range("B121").Value = 2734
$ws.Range("C121").Value = 45
$ws.Range("D121").Value = 2010
$ws.Range("E121").Value = 697

# Row 139
$ws.Range("B139").Value = 1483
$ws.Range("C139").Value = 88
$ws.Range("D139").Value = 520
$ws.Range("E139").Value = 899
$ws.Range("G139").Value = 2
$ws.Range("H139").Value = 64

# Row 149
$ws.Range("B149").Value = 1012
$ws.Range("C149").Value = 11
$ws.Range("D149").Value = 697
$ws.Range("E149").Value = 293
$ws.Range("G149").Value = 1
$ws.Range("H149").Value = 22

$wb.Save()
